$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update BANKNIFTY and FINNIFTY expiry dates (shift forward by one week)
$ws.Range("B6").Value = 45441
$ws.Range("B7").Value = 45448
$ws.Range("B8").Value = 45440
$ws.Range("B9").Value = 45447
